$d = $word.ActiveDocument
$bullet = [char]0x2022

# --- 1. Collapse the three CORE COMPETENCIES detail paragraphs into one condensed line ---
# Locate the three paragraphs by their distinctive leading text instead of a hard-coded
# paragraph index, so the script is resilient to any small structural differences.
$surveyPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Survey Methodology & Research Design:")) {
        $surveyPara = $p
        break
    }
}

if ($surveyPara -eq $null) {
    throw "Could not locate the 'Survey Methodology & Research Design:' paragraph"
}

$surveyIndex = $surveyPara.Index
# The next two paragraphs are "Redistricting & Geospatial Analysis: ..." and
# "Data Analysis & Visualization: ...". Remove them (higher index first) and then
# overwrite the first paragraph's text with the condensed summary line.
$d.Paragraphs.Item($surveyIndex + 2).Range.Delete()
$d.Paragraphs.Item($surveyIndex + 1).Range.Delete()

$coreCompText = "Survey Methodology & Research Design $bullet Redistricting & Geospatial Analysis $bullet Data Analysis & Visualization"
$d.Paragraphs.Item($surveyIndex).Range.Text = $coreCompText

# --- 2. Append a new "TECHNICAL SKILLS" section just before the closing paragraph ---
# Anchor on the last bullet of the "Survey Methodology & Data Analysis" achievements
# block ("Developed advanced segmentation models...") which immediately precedes the
# document's final "For a more detailed..." paragraph.
$segmentationPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("$bullet Developed advanced segmentation models")) {
        $segmentationPara = $p
    }
}

if ($segmentationPara -eq $null) {
    throw "Could not locate the 'Developed advanced segmentation models' paragraph"
}

$segmentationPara.Range.InsertParagraphAfter()
$headingPara = $d.Paragraphs.Item($segmentationPara.Index + 1)
$headingPara.Range.Text = "TECHNICAL SKILLS"
$headingPara.Style = "Heading 2"

$detailTexts = @(
    "SURVEY METHODOLOGY & RESEARCH DESIGN Survey Design and Questionnaire Development for Political and Market Research; Sampling Methodology and Statistical Analysis (R, SPSS, Stata, OSCAR); Random Device Engagement (RDE), Text Message, Web Panel, and Live Telephone Calling; Expert Testimony and Consultation on Research Methodology",
    "REDISTRICTING & GEOSPATIAL ANALYSIS Redistricting Software Development and Boundary Estimation Systems; Geospatial Analysis; Choropleths and Hexagonal Grid Maps for Demographic Visualization; Court Case Analysis and Expert Testimony for Redistricting",
    "DATA ANALYSIS & VISUALIZATION Advanced Statistical Modeling and Analysis (Regression, Clustering, Segmentation); Data Visualization; Consumer Behavior Analysis and Market Segmentation; Multi-million Dollar Research Project Management"
)

$anchor = $headingPara
foreach ($t in $detailTexts) {
    $anchor.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($anchor.Index + 1)
    $newPara.Style = "Normal"
    $newPara.Range.Text = $t
    $anchor = $newPara
}
